$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124, shifting existing rows 124:139 down to 125:140
$ws.Rows("124:124").Insert()

# Populate the newly inserted row 124 with the new data record
$ws.Range("A124").Value = 1
$ws.Range("B124").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C124").Value = "Arica y Parinacota"
$ws.Range("D124").Value = 44951
$ws.Range("D124").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E124").Value = 15
$ws.Range("F124").Value = "Fruta"
$ws.Range("G124").Value = 100102
$ws.Range("H124").Value = "Cítricos"
$ws.Range("I124").Value = 100102004
$ws.Range("J124").Value = "Mandarina"
$ws.Range("K124").Value = "Murcott"
$ws.Range("L124").Value = "Segunda"
$ws.Range("M124").Value = 550
$ws.Range("N124").Value = 17000
$ws.Range("O124").Value = 18000
$ws.Range("P124").Value = 17455
$ws.Range("Q124").Value = "`$/caja 20 kilos"
$ws.Range("R124").Value = "Región de Coquimbo"
$ws.Range("S124").Value = 873
$ws.Range("T124").Value = 20
